# Update "想去人数" (want-to-go count) figures in the 广州-漫展信息 workbook.
# Sheet "展览" (index 1), "演出" (index 2) and "全部类型" (index 4) each get
# updated F-column values. "全部类型" aggregates rows from the other sheets,
# so several of the same logical updates appear twice (once in the source
# sheet, once in the aggregated sheet).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibition) sheet ---
$wsExhibit.Range("F6").Value  = 942
$wsExhibit.Range("F8").Value  = 17
$wsExhibit.Range("F9").Value  = 984
$wsExhibit.Range("F10").Value = 781
$wsExhibit.Range("F16").Value = 570
$wsExhibit.Range("F21").Value = 1140
$wsExhibit.Range("F23").Value = 1361
$wsExhibit.Range("F24").Value = 674
$wsExhibit.Range("F26").Value = 1256
$wsExhibit.Range("F30").Value = 2123
$wsExhibit.Range("F31").Value = 194
$wsExhibit.Range("F32").Value = 158
$wsExhibit.Range("F33").Value = 1367

# --- 演出 (Show) sheet ---
$wsShow.Range("F3").Value = 516

# --- 全部类型 (All types) sheet ---
$wsAll.Range("F8").Value  = 516
$wsAll.Range("F12").Value = 942
$wsAll.Range("F15").Value = 17
$wsAll.Range("F16").Value = 984
$wsAll.Range("F17").Value = 781
$wsAll.Range("F28").Value = 570
$wsAll.Range("F33").Value = 1140
$wsAll.Range("F35").Value = 1361
$wsAll.Range("F36").Value = 674
$wsAll.Range("F38").Value = 1256
$wsAll.Range("F44").Value = 2123
$wsAll.Range("F45").Value = 194
$wsAll.Range("F46").Value = 158
$wsAll.Range("F47").Value = 1367
